$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I (2015 data): set values first ---
$ws.Range("I1").Value = 2015
$ws.Range("I2").Value = 28.479772567749023
$ws.Range("I3").Value = 25.384542465209961
$ws.Range("I4").Value = 1.0357397794723511
$ws.Range("I5").Value = 0.94090276956558228
$ws.Range("I6").Value = 35.414726257324219
$ws.Range("I7").Value = 29.453279495239258
$ws.Range("I8").Value = 10.538521766662598
$ws.Range("I9").Value = 9.088623046875
$ws.Range("I10").Value = 14.993505477905273
$ws.Range("I11").Value = 11.510863304138184
$ws.Range("I12").Value = 9.4726734161376953
$ws.Range("I13").Value = 7.7272524833679199
$ws.Range("I14").Value = 16.91871452331543
$ws.Range("I15").Value = 2.1791555881500244
$ws.Range("I16").Value = 21.804689407348633
$ws.Range("I17").Value = 18.909276962280273
$ws.Range("I18").Value = 11.796778678894043
$ws.Range("I19").Value = 2.3814046382904053
$ws.Range("I20").Value = 12.565120697021484
$ws.Range("I21").Value = 12.479278564453125
$ws.Range("I22").Value = 0.96558243036270142
$ws.Range("I23").Value = 23.983844757080078
$ws.Range("I24").Value = 30.7060546875
$ws.Range("I25").Value = 14.946524620056152
$ws.Range("I26").Value = 13.444860458374023
$ws.Range("I27").Value = 20.402881622314453
$ws.Range("I28").Value = 70.648017883300781
$ws.Range("I32").Value = 20.402881622314453
$ws.Range("I33").Value = 8.3964433670043945
$ws.Range("I34").Value = 3.4136171340942383
$ws.Range("I35").Value = 4.0844755172729492

# --- New empty-but-styled column I cells (rows 29-31) ---
# I29 left blank, style applied below
# I30 left blank, style applied below
# I31 left blank, style applied below

# --- Apply number format to the whole new column I so style matches existing columns (s=1) ---
$ws.Range("I1:I35").NumberFormat = "0"

# --- Updated values for rows 14-26, columns B:H (existing columns) ---
$ws.Range("B14").Value = 16.91871452331543
$ws.Range("C14").Value = 16.91871452331543
$ws.Range("D14").Value = 16.91871452331543
$ws.Range("E14").Value = 16.91871452331543
$ws.Range("F14").Value = 16.91871452331543
$ws.Range("G14").Value = 16.91871452331543
$ws.Range("H14").Value = 16.91871452331543
$ws.Range("B15").Value = 2.2051389217376709
$ws.Range("C15").Value = 2.2049911022186279
$ws.Range("D15").Value = 2.2061336040496826
$ws.Range("E15").Value = 2.1997659206390381
$ws.Range("F15").Value = 2.1909668445587158
$ws.Range("G15").Value = 2.1853375434875488
$ws.Range("H15").Value = 2.1803452968597412
$ws.Range("B16").Value = 21.881288528442383
$ws.Range("C16").Value = 21.865360260009766
$ws.Range("D16").Value = 21.846712112426758
$ws.Range("E16").Value = 21.840641021728516
$ws.Range("F16").Value = 21.825418472290039
$ws.Range("G16").Value = 21.806716918945313
$ws.Range("H16").Value = 21.805488586425781
$ws.Range("B17").Value = 18.838916778564453
$ws.Range("C17").Value = 18.850732803344727
$ws.Range("D17").Value = 18.862716674804688
$ws.Range("E17").Value = 18.876714706420898
$ws.Range("F17").Value = 18.888875961303711
$ws.Range("G17").Value = 18.904705047607422
$ws.Range("H17").Value = 18.908233642578125
$ws.Range("B18").Value = 11.785534858703613
$ws.Range("C18").Value = 11.789078712463379
$ws.Range("D18").Value = 11.793749809265137
$ws.Range("E18").Value = 11.792673110961914
$ws.Range("F18").Value = 11.792230606079102
$ws.Range("G18").Value = 11.791597366333008
$ws.Range("H18").Value = 11.794395446777344
$ws.Range("B19").Value = 2.3810474872589111
$ws.Range("C19").Value = 2.3810431957244873
$ws.Range("D19").Value = 2.3811821937561035
$ws.Range("E19").Value = 2.3813536167144775
$ws.Range("F19").Value = 2.3815581798553467
$ws.Range("G19").Value = 2.3813114166259766
$ws.Range("H19").Value = 2.3813951015472412
$ws.Range("B20").Value = 12.550850868225098
$ws.Range("C20").Value = 12.552334785461426
$ws.Range("D20").Value = 12.553961753845215
$ws.Range("E20").Value = 12.551446914672852
$ws.Range("F20").Value = 12.560141563415527
$ws.Range("G20").Value = 12.568150520324707
$ws.Range("H20").Value = 12.567581176757813
$ws.Range("B21").Value = 12.473625183105469
$ws.Range("C21").Value = 12.472444534301758
$ws.Range("D21").Value = 12.471128463745117
$ws.Range("E21").Value = 12.472867012023926
$ws.Range("F21").Value = 12.476411819458008
$ws.Range("G21").Value = 12.477934837341309
$ws.Range("H21").Value = 12.478243827819824
$ws.Range("B22").Value = 0.96488356590270996
$ws.Range("C22").Value = 0.96529960632324219
$ws.Range("D22").Value = 0.96570080518722534
$ws.Range("E22").Value = 0.96582388877868652
$ws.Range("F22").Value = 0.96568220853805542
$ws.Range("G22").Value = 0.9655308723449707
$ws.Range("H22").Value = 0.96560186147689819
$ws.Range("B23").Value = 24.086427688598633
$ws.Range("C23").Value = 24.070352554321289
$ws.Range("D23").Value = 24.052845001220703
$ws.Range("E23").Value = 24.040407180786133
$ws.Range("F23").Value = 24.016386032104492
$ws.Range("G23").Value = 23.992055892944336
$ws.Range("H23").Value = 23.985834121704102
$ws.Range("B24").Value = 30.62445068359375
$ws.Range("C24").Value = 30.639812469482422
$ws.Range("D24").Value = 30.656467437744141
$ws.Range("E24").Value = 30.669387817382813
$ws.Range("F24").Value = 30.681106567382813
$ws.Range("G24").Value = 30.69630241394043
$ws.Range("H24").Value = 30.702629089355469
$ws.Range("B25").Value = 14.93189811706543
$ws.Range("C25").Value = 14.933378219604492
$ws.Range("D25").Value = 14.935144424438477
$ws.Range("E25").Value = 14.932801246643066
$ws.Range("F25").Value = 14.941699981689453
$ws.Range("G25").Value = 14.949461936950684
$ws.Range("H25").Value = 14.948976516723633
$ws.Range("B26").Value = 13.438508033752441
$ws.Range("C26").Value = 13.437744140625
$ws.Range("D26").Value = 13.436829566955566
$ws.Range("E26").Value = 13.438690185546875
$ws.Range("F26").Value = 13.442093849182129
$ws.Range("G26").Value = 13.443466186523438
$ws.Range("H26").Value = 13.443845748901367
